$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily mark the Price column as text so that purely numeric-looking
# values (e.g. "248.43") are stored as strings, matching the inlineStr cells
# already in the workbook, instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '37.171.35'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.052.82'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '248.43'
$ws.Range("E5").Value = '  -1.78%  '
$ws.Range("D6").Value = '0.666'
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '57.03'
$ws.Range("E8").Value = '  -1.13%  '
$ws.Range("E9").Value = '  -0.33%  '
$ws.Range("E10").Value = '  -2.23%  '
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = '16.25'
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("D13").Value = '0.922'
$ws.Range("E13").Value = '  +14.07%  '
$ws.Range("D14").Value = '2.348.40'
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").Value = '5.76'
$ws.Range("E15").Value = '  +2.47%  '
$ws.Range("D16").Value = '2.051.01'
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("D17").Value = '18.74'
$ws.Range("E17").Value = '  +13.27%  '
$ws.Range("D18").Value = '37.189.85'
$ws.Range("E18").Value = '  +0.39%  '
$ws.Range("D19").Value = '74.64'
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("D20").Value = [string]::Concat('0.0', [char]0x2083, '0899')
$ws.Range("E20").Value = '  -2.49%  '
$ws.Range("D21").Value = '5.47'
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").Value = '237.33'
$ws.Range("E22").Value = '  -0.47%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("E24").Value = '  +3.42%  '
$ws.Range("D25").Value = '9.67'
$ws.Range("E25").Value = '  +3.78%  '
$ws.Range("D26").Value = '170.16'
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("D27").Value = '2.18'
$ws.Range("E27").Value = '  -4.42%  '
$ws.Range("D28").Value = '20.23'
$ws.Range("E28").Value = '  -0.32%  '
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("D30").Value = '5.05'
$ws.Range("E30").Value = '  +6.00%  '
$ws.Range("E31").Value = '  +1.53%  '
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("D33").Value = '4.55'
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("D34").Value = '0.0881'
$ws.Range("E34").Value = '  -1.91%  '
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("D37").Value = '1.78'
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("D39").Value = '5.32'
$ws.Range("E39").Value = '  +13.49%  '
$ws.Range("D40").Value = '3.09'
$ws.Range("E40").Value = '  +7.97%  '
$ws.Range("D41").Value = '0.0988'
$ws.Range("E41").Value = '  -13.06%  '
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("D43").Value = '17.66'
$ws.Range("E43").Value = '  -1.01%  '
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("D45").Value = '96.38'
$ws.Range("E45").Value = '  -1.33%  '
$ws.Range("E46").Value = '  -1.20%  '
$ws.Range("D47").Value = '1.274.45'
$ws.Range("E47").Value = '  -1.47%  '
$ws.Range("E48").Value = '  -2.03%  '
$ws.Range("D49").Value = '6.84'
$ws.Range("E49").Value = '  -1.15%  '
$ws.Range("D50").Value = '2.237.03'
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("D51").Value = '44.52'
$ws.Range("E51").Value = '  +1.01%  '

# Restore the column to the default (General) style now that the text values
# are committed, so untouched/updated cells keep style index 0 like the source.
$ws.Range("D2:D51").Style = "Normal"
